$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$ws = $wb.Worksheets.Item("Metadata")

# Version: 6.1.1 -> 7.0.0
$ws.Range("B3").Value = "7.0.0"

# Title: "Snapshot Age in months" -> "Snapshot Age in Months"
$ws.Range("B5").Value = "Snapshot Age in Months"

# Date: 2022-06-06T15:56:40+00:00 -> 2022-09-01T20:48:10+00:00
$ws.Range("B8").Value = "2022-09-01T20:48:10+00:00"
